# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets
# of the workbook, per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 652
$ws1.Range("G3").Value  = 70
$ws1.Range("F4").Value  = 813
$ws1.Range("G4").Value  = 70
$ws1.Range("F5").Value  = 498
$ws1.Range("F8").Value  = 886
$ws1.Range("F10").Value = 819
$ws1.Range("F11").Value = 644
$ws1.Range("F12").Value = 118
$ws1.Range("F15").Value = 722
$ws1.Range("F16").Value = 212
$ws1.Range("F18").Value = 462
$ws1.Range("F19").Value = 1228
$ws1.Range("F21").Value = 937
$ws1.Range("F22").Value = 2682
$ws1.Range("F23").Value = 1154
$ws1.Range("F24").Value = 612
$ws1.Range("F25").Value = 139
$ws1.Range("F26").Value = 1197
$ws1.Range("F28").Value = 893
$ws1.Range("F29").Value = 93
$ws1.Range("F30").Value = 1226

# --- Sheet 2: 演出 ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F10").Value = 20
$ws2.Range("F12").Value = 17

# --- Sheet 3: 本地生活 -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 698

# --- Sheet 4: 全部类型 -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 698
$ws4.Range("F4").Value  = 652
$ws4.Range("G4").Value  = 70
$ws4.Range("F5").Value  = 813
$ws4.Range("G5").Value  = 70
$ws4.Range("F6").Value  = 498
$ws4.Range("F14").Value = 886
$ws4.Range("F16").Value = 820
$ws4.Range("F17").Value = 644
$ws4.Range("F18").Value = 118
$ws4.Range("F23").Value = 20
$ws4.Range("F26").Value = 722
$ws4.Range("F27").Value = 212
$ws4.Range("F29").Value = 462
$ws4.Range("F30").Value = 1228
$ws4.Range("F32").Value = 937
$ws4.Range("F33").Value = 2682
$ws4.Range("F34").Value = 1154
$ws4.Range("F35").Value = 612
$ws4.Range("F36").Value = 139
$ws4.Range("F37").Value = 1197
$ws4.Range("F39").Value = 17
$ws4.Range("F40").Value = 893
$ws4.Range("F41").Value = 93
$ws4.Range("F42").Value = 1226
